# Adds three new lineage rows (first_name, full_name, last_name) for the
# "report_table" target just before the existing trailing "orders" row,
# and fixes up the upstream-lineage formula text for the
# customer_orders_table.account_length_days row (E9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the upstream lineage text on row 9 (account_length_days).
$ws.Range("E9").Value = "datediff(current_date(), CAST(account_open_date AS DATE))`nfirst(account_length_days)"

# 2) The sheet currently ends at row 33 with the "report_table" / "orders"
#    lineage row. Push that row down by three and insert new rows for
#    first_name, full_name and last_name above it (same target
#    catalog/database/table as the existing "orders" row).
$catalog = $ws.Range("A33").Text
$database = $ws.Range("B33").Text
$table = $ws.Range("C33").Text

# Copy the existing data-row formatting (style index 2: wrap text, vertical
# centered) down onto the three new rows plus the row that will hold the
# relocated "orders" entry.
$ws.Range("A33:E33").Copy()
$ws.Range("A34:E36").PasteSpecial(-4122)

$rows = @(
    @("first_name", "first_name"),
    @("full_name", "full_name"),
    @("last_name", "last_name"),
    @("orders", "orders")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = 33 + $i
    $colD = $rows[$i][0]
    $colE = $rows[$i][1]

    $ws.Range("A" + $r).Value = $catalog
    $ws.Range("B" + $r).Value = $database
    $ws.Range("C" + $r).Value = $table
    $ws.Range("D" + $r).Value = $colD
    $ws.Range("E" + $r).Value = $colE
}
